# Update media and conferences:
#  - London / MRF AMR Annual Conference 2021 is now "Attended" (row 5,
#    previously "Will attend").
#  - A new row is inserted for London / MRF AMR Annual Conference 2021 /
#    Aug-21 / Will attend, pushing the Bologna / EPIDEMICS8 row down from
#    row 7 to row 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Microbiology Society Annual Conference 2021, London) - now attended.
$ws.Range("D5").Value = "Attended"

# Insert a new row before the current row 7 (Bologna/EPIDEMICS8), shifting
# it down to row 8 and carrying its formatting (style) along.
$ws.Rows.Item(7).Insert()

# Fill in the new row 7 with the MRF AMR Annual Conference 2021 details.
$ws.Range("A7").Value = "London"
$ws.Range("C7").Value = "Aug-21"
$ws.Range("B7").Value = "MRF AMR Annual Conference 2021"
$ws.Range("D7").Value = "Will attend"
$ws.Range("E7").Value = 51.609864999999999
$ws.Range("F7").Value = -0.21809200000000001

# Match the recorded selection after the edit.
$ws.Range("G8").Select()
